$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.641980290412903
$ws.Range("B1").Value = 2.887482166290283
$ws.Range("C1").Value = 3.715534210205078
$ws.Range("D1").Value = 1.443675398826599
$ws.Range("E1").Value = 0.9655657410621643
